$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 68 (last existing data row) ---
# The "Ordered Amount" B value was mis-keyed; correct it.
$ws.Range("B68").Value = 75920
# The "Ordered Amount" label cell had picked up a stray italic font; make it
# match the rest of that column (D63:D66).
$ws.Range("D68").Font.Italic = $false
$ws.Range("E68").Formula = '=IF(A68="","",SUM(E67-B68+C68))'

# --- New ledger rows, entered through 22 Feb 8AM per the commit message ---

# Row 69: 20 Feb 2020, Amt Debited 34320, Ordered Amount
$ws.Range("A69").Value = 43881
$ws.Range("A69").NumberFormat = "[$-409]d\-mmm\-yyyy;@"
$ws.Range("B69").Value = 34320
$ws.Range("D69").Value = "Ordered Amount"
$ws.Range("D69").Font.Italic = $false
$ws.Range("E69").Formula = '=IF(A69="","",SUM(E68-B69+C69))'

# Row 70: 21 Feb 2020, Amt Credited 147826, Manual Added
$ws.Range("A70").Value = 43882
$ws.Range("A70").NumberFormat = "[$-409]d\-mmm\-yyyy;@"
$ws.Range("C70").Value = 147826
$ws.Range("D70").Value = "Manual Added"
$ws.Range("E70").Formula = '=IF(A70="","",SUM(E69-B70+C70))'

# Row 71: 21 Feb 2020, Amt Debited 29120, Ordered Amount
$ws.Range("A71").Value = 43882
$ws.Range("A71").NumberFormat = "[$-409]d\-mmm\-yyyy;@"
$ws.Range("B71").Value = 29120
$ws.Range("D71").Value = "Ordered Amount"
$ws.Range("D71").Font.Italic = $false
$ws.Range("E71").Formula = '=IF(A71="","",SUM(E70-B71+C71))'

# Move the frozen-pane top-left cell and selection down one row to track the
# newly entered data, same as the author's cursor position when they saved.
$ws.Range("D68").Select()
